# Apply cryptos list price/volume refresh (GitHub Actions data update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.670.23'
$ws.Range('E2').Value = '  -0.28%  '

$ws.Range('D3').Value = '3.730.24'
$ws.Range('E3').Value = '  -1.96%  '

$ws.Range('D4').Value = '''1.00'

$ws.Range('D5').Value = '''591.37'
$ws.Range('E5').Value = '  -1.30%  '

$ws.Range('D6').Value = '''165.59'
$ws.Range('E6').Value = '  -1.90%  '

$ws.Range('D7').Value = '3.727.99'
$ws.Range('E7').Value = '  -1.98%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('E9').Value = '  -2.04%  '

$ws.Range('D10').Value = '''0.159'
$ws.Range('E10').Value = '  -3.43%  '

$ws.Range('D11').Value = '''6.49'
$ws.Range('E11').Value = '  -0.19%  '

$ws.Range('D12').Value = '''0.450'
$ws.Range('E12').Value = '  -2.61%  '

$ws.Range('D13').Value = '''0.0000261'
$ws.Range('E13').Value = '  -4.92%  '

$ws.Range('D14').Value = '''36.09'
$ws.Range('E14').Value = '  -2.14%  '

$ws.Range('D15').Value = '4.352.43'
$ws.Range('E15').Value = '  -2.18%  '

$ws.Range('D16').Value = '3.731.06'
$ws.Range('E16').Value = '  -1.74%  '

$ws.Range('D17').Value = '67.630.53'
$ws.Range('E17').Value = '  -0.53%  '

$ws.Range('D18').Value = '''18.30'
$ws.Range('E18').Value = '  -0.92%  '

$ws.Range('D19').Value = '''7.02'
$ws.Range('E19').Value = '  -5.27%  '

$ws.Range('E20').Value = '  -0.34%  '

$ws.Range('D21').Value = '''10.72'
$ws.Range('E21').Value = '  -1.41%  '

$ws.Range('D22').Value = '''467.61'
$ws.Range('E22').Value = '  -0.34%  '

$ws.Range('E23').Value = '  -5.04%  '

$ws.Range('D24').Value = '''82.73'
$ws.Range('E24').Value = '  -0.70%  '

$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').Value = '''2.18'
$ws.Range('E25').Value = '  -6.21%  '

$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '''0.0000133'
$ws.Range('E26').Value = '  -12.04%  '

$ws.Range('D27').Value = '''11.99'
$ws.Range('E27').Value = '  -1.84%  '

$ws.Range('D28').Value = '''10.08'
$ws.Range('E28').Value = '  -2.10%  '

$ws.Range('E29').Value = '  +0.00%  '

$ws.Range('D30').Value = '3.872.45'

$ws.Range('E31').Value = '  -5.60%  '

$ws.Range('D32').Value = '''7.32'
$ws.Range('E32').Value = '  -5.42%  '

$ws.Range('E33').Value = '  -2.83%  '

$ws.Range('D34').Value = '''29.50'
$ws.Range('E34').Value = '  -4.04%  '

$ws.Range('D35').Value = '''9.03'
$ws.Range('E35').Value = '  -3.21%  '

$ws.Range('D36').Value = '3.679.61'
$ws.Range('E36').Value = '  -2.47%  '

$ws.Range('E37').Value = '  -5.41%  '

$ws.Range('D38').Value = '''3.41'
$ws.Range('E38').Value = '  -9.14%  '

$ws.Range('D39').Value = '''0.138'
$ws.Range('E39').Value = '  -1.44%  '

$ws.Range('D40').Value = '''0.989'
$ws.Range('E40').Value = '  -1.78%  '

$ws.Range('D41').Value = '''5.74'
$ws.Range('E41').Value = '  -4.22%  '

$ws.Range('E42').Value = '  -0.11%  '

$ws.Range('D44').Value = '''0.305'
$ws.Range('E44').Value = '  -4.76%  '

$ws.Range('D45').Value = '''8.52'
$ws.Range('E45').Value = '  -3.10%  '

$ws.Range('D46').Value = '''1.91'
$ws.Range('E46').Value = '  -3.39%  '

$ws.Range('D47').Value = '''45.22'
$ws.Range('E47').Value = '  -2.55%  '

$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '''387.05'
$ws.Range('E48').Value = '  -5.86%  '

$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '''142.74'
$ws.Range('E49').Value = '  +0.38%  '

$ws.Range('D50').Value = '''25.41'
$ws.Range('E50').Value = '  +0.35%  '

$ws.Range('D51').Value = '''0.0346'
$ws.Range('E51').Value = '  -3.97%  '
